# 823-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-Late Repayment-Makerepayment1.xlsx
# Apply the "Late Repayment" re-run: updated repayment figures after a late
# payment recompute, plus the resulting sheet/selection navigation and a
# header-formatting cleanup on the Transactions sheet.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsRepay   = $wb.Worksheets.Item("Repayment Schedule")
$wsTrans   = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------
# Summary sheet: recomputed principal / balance figures
# ---------------------------------------------------------------------
$wsSummary.Range("A3").Value = 675.05
$wsSummary.Range("E3").Value = 525.05

# ---------------------------------------------------------------------
# Repayment Schedule sheet: recomputed late-repayment schedule (rows 7-15)
# ---------------------------------------------------------------------
$wsRepay.Range("F7").Value  = 812.72
$wsRepay.Range("G7").Value  = 6764.43
$wsRepay.Range("H7").Value  = 75.77

$wsRepay.Range("F8").Value  = 820.85
$wsRepay.Range("G8").Value  = 5943.58
$wsRepay.Range("H8").Value  = 67.64

$wsRepay.Range("F9").Value  = 829.05
$wsRepay.Range("G9").Value  = 5114.53
$wsRepay.Range("H9").Value  = 59.44

$wsRepay.Range("F10").Value = 837.34
$wsRepay.Range("G10").Value = 4277.19
$wsRepay.Range("H10").Value = 51.15

$wsRepay.Range("F11").Value = 845.72
$wsRepay.Range("G11").Value = 3431.47
$wsRepay.Range("H11").Value = 42.77

$wsRepay.Range("F12").Value = 854.18
$wsRepay.Range("G12").Value = 2577.29
$wsRepay.Range("H12").Value = 34.31

$wsRepay.Range("F13").Value = 862.72
$wsRepay.Range("G13").Value = 1714.57
$wsRepay.Range("H13").Value = 25.77

$wsRepay.Range("F14").Value = 871.34
$wsRepay.Range("G14").Value = 843.23
$wsRepay.Range("H14").Value = 17.15

$wsRepay.Range("F15").Value = 843.23
$wsRepay.Range("H15").Value = 8.43
$wsRepay.Range("K15").Value = 851.66
$wsRepay.Range("P15").Value = 851.66

# ---------------------------------------------------------------------
# Transactions sheet: transaction ids recomputed + header row formatting
# unified to a single consistent style (same as already used on C1:J1)
# ---------------------------------------------------------------------
$wsTrans.Range("A2").Value = 183
$wsTrans.Range("A3").Value = 139

$wsTrans.Range("C1").Copy() | Out-Null
$wsTrans.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping, matching the saved view state
# ---------------------------------------------------------------------
$wsSummary.Range("C4").Select() | Out-Null
$wsRepay.Range("D10").Select() | Out-Null

$wsTrans.Activate() | Out-Null
$wsTrans.Range("I3").Select() | Out-Null
